# final version of the report (delivered)
#
# Applies the tracked-changes style rewording of the "Actor Architecture"
# report: a few word-level tweaks, a couple of paragraph merges/splits,
# a fair amount of new prose in the Add/Contains/Remove message
# paragraphs, and a uniform single-line-spacing (w:spacing
# w:line="240" w:lineRule="auto") applied to every paragraph from
# "Actor Architecture:" onward.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "The actor architecture is used ..." paragraph (#4): small tweak.
# ---------------------------------------------------------------------
$r = $d.Paragraphs(4).Range
$r.Find.Execute("but can only affect and communicate", $true, $false, $false, $false, $false, $true, 1, $false, "but they can affect and communicate", 2)

# ---------------------------------------------------------------------
# 2. Drop the blank paragraph between that paragraph and "Implementation:"
#    (#5 becomes absorbed, "Implementation:" shifts up to #5).
# ---------------------------------------------------------------------
$d.Paragraphs(5).Range.Delete()

# ---------------------------------------------------------------------
# 3. "To implement a Binary Tree ..." paragraph (now #6): actors -> Actors
# ---------------------------------------------------------------------
$r = $d.Paragraphs(6).Range
$r.Find.Execute("uses actors the actors have", $true, $false, $false, $false, $false, $true, 1, $false, "uses Actors the Actors have", 2)

# ---------------------------------------------------------------------
# 4. "A left and right pointer to other actors" (now #9)
# ---------------------------------------------------------------------
$r = $d.Paragraphs(9).Range
$r.Find.Execute("A left and right pointer to other actors", $true, $false, $false, $false, $false, $true, 1, $false, "A left and right reference to other Actors", 2)

# ---------------------------------------------------------------------
# 5. "Each Actor has its own thread ..." paragraph (now #10):
#    "leaf node corresponds to." -> "leaf Actor corresponds to."
# ---------------------------------------------------------------------
$r = $d.Paragraphs(10).Range
$r.Find.Execute("the side of the leaf node corresponds to.", $true, $false, $false, $false, $false, $true, 1, $false, "the side of the leaf Actor corresponds to.", 2)

# ---------------------------------------------------------------------
# 6. Split paragraph (now #11) into two: the "To find the recipient..."
#    sentence stays, and a new paragraph starts at "The add message...".
# ---------------------------------------------------------------------
$p11 = $d.Paragraphs(11).Range
$splitPoint = $d.Range($p11.Start, $p11.End)
$splitPoint.Find.Execute("The add message searches")
$splitPoint.Collapse(1)
$splitPoint.InsertParagraphBefore()
# drop the trailing space left at the end of the first half
$r = $d.Paragraphs(11).Range
$r.Find.Execute("it being superior. ", $true, $false, $false, $false, $false, $true, 1, $false, "it being superior.", 2)

# ---------------------------------------------------------------------
# 7. Rewrite the (now split-off) "The add message ..." paragraph (#12)
#    with the final, expanded wording.
# ---------------------------------------------------------------------
$addOld = "The add message searches the tree to find if the number in the message is already present in the tree, in case it reaches the bottom of the tree without finding an Actor with the number of the message it then is created an Actor with that number, and the Actor responsible for creating the new Actor sends a AddResponse message to the RootActor."
$addNew = "The add message is forwarded down the tree to find if the number in the message is already present in the tree, in case it reaches the bottom of the tree without finding an Actor with the number of the message it is then created an Actor with that number. The Actor responsible for creating the new Actor sends an AddResponse message to the RootActor. In the case the number is already attributed to an Actor present in the tree it is sent an AddResponse message to the RootActor notifying the number is already present."
$r = $d.Paragraphs(12).Range
$r.Find.Execute($addOld, $true, $false, $false, $false, $false, $true, 1, $false, $addNew, 2)

# ---------------------------------------------------------------------
# 8. Rewrite "The contains message ..." paragraph (#13).
# ---------------------------------------------------------------------
$containsOld = "The contains message searches the tree to find the number in the message, when it finds the Actor with the number present in the message or it reaches the bottom of the tree without finding the number it sends a ContainsResponse message to the RootActor."
$containsNew = "The contains message is forwarded down the tree to find if the number in the message is present matches an Actor present in the tree. When it finds the Actor with the number present in the message or it reaches the bottom of the tree without finding the number it sends a ContainsResponse message to the RootActor."
$r = $d.Paragraphs(13).Range
$r.Find.Execute($containsOld, $true, $false, $false, $false, $false, $true, 1, $false, $containsNew, 2)

# ---------------------------------------------------------------------
# 9. Expand "The remove message" paragraph (#14) with the full
#    description (the _GoBack bookmark that used to sit at the end of
#    this paragraph is dropped along with the old short text).
# ---------------------------------------------------------------------
$removeNew = "The remove message is forwarded down the tree till it finds the Actor with the number present in the message. When that Actor is found, it then proceeds to determine how many leaf Actors it has. In the case of only one it sends a message to his supervisor with the reference of that Actor. In the case it has two leaf Actors it starts a process to find the Actor with the smallest number on its right subtree. Once it’s found that Actor sends a message to the supervisor of the Actor that it was removed with his reference so that it can replace the position of the removed Actor."
$r = $d.Paragraphs(14).Range
$r.Text = $removeNew

# ---------------------------------------------------------------------
# 10. Insert a brand-new paragraph after it ("Regardless of the number
#     of leaf Actors ..."), reusing the two blank paragraphs (#15, #16)
#     that used to separate "The remove message" from "Garbage
#     Collection:" — fill the first with the new text, delete the
#     second so only one paragraph remains.
# ---------------------------------------------------------------------
$regardlessText = "Regardless of the number of leaf Actors the removed Actor then deletes all messages present on his mailbox, the references it contained of his leaf Actors and signals himself to stop processing messages, followed by sending a RemoveResponse message to the RootActor signalling the success of the operation."
$d.Paragraphs(15).Range.Text = $regardlessText
$d.Paragraphs(16).Range.Delete()

# ---------------------------------------------------------------------
# 11. Apply single line spacing (w:spacing w:line="240"
#     w:lineRule="auto") to every paragraph from "Actor Architecture:"
#     onward.
# ---------------------------------------------------------------------
for ($i = 3; $i -le $d.Paragraphs.Count; $i++) {
    $d.Paragraphs($i).Range.ParagraphFormat.LineSpacingRule = 0
}
